$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 751.7222
$ws.Range("I2").Value2 = 266.44446
$ws.Range("J2").Value2 = 1237
$ws.Range("K2").Value2 = 266.44446
$ws.Range("L2").Value2 = 1237
$ws.Range("M2").Value2 = -153.44446
$ws.Range("N2").Value2 = -1463

$ws.Range("H4").Value2 = 1996.6666
$ws.Range("I4").Value2 = 0
$ws.Range("J4").Value2 = 1996.6666
$ws.Range("K4").Value2 = 0
$ws.Range("L4").Value2 = 1996.6666
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value2 = -2224.6666

$ws.Range("H8").Value2 = 1086.4546
$ws.Range("I8").Value2 = 26.142857
$ws.Range("J8").Value2 = 2942
$ws.Range("K8").Value2 = 78.42857100000001
$ws.Range("L8").Value2 = 8826
$ws.Range("M8").Value2 = 60.57142899999999
$ws.Range("N8").Value2 = -9104

$ws.Range("H9").Value2 = 993.9375
$ws.Range("I9").Value2 = 1600
$ws.Range("J9").Value2 = 387.875
$ws.Range("K9").Value2 = 1600
$ws.Range("L9").Value2 = 387.875
$ws.Range("M9").Value2 = -1431
$ws.Range("N9").Value2 = -725.875

$ws.Range("H32").Value2 = 2703.3076
$ws.Range("I32").Value2 = 872.75
$ws.Range("J32").Value2 = 3516.889
$ws.Range("K32").Value2 = 872.75
$ws.Range("L32").Value2 = 3516.889
$ws.Range("M32").Value2 = -546.75
$ws.Range("N32").Value2 = -4168.889

$ws.Range("H38").Value2 = 9093813
$ws.Range("I38").Value2 = 11112104
$ws.Range("J38").Value2 = 11500
$ws.Range("K38").Value2 = 33336312
$ws.Range("L38").Value2 = 34500
$ws.Range("M38").Value2 = -33335940
$ws.Range("N38").Value2 = -35244

$ws.Range("H40").Value2 = 1749.6923
$ws.Range("I40").Value2 = 1470.4117
$ws.Range("J40").Value2 = 2277.2222
$ws.Range("K40").Value2 = 1470.4117
$ws.Range("L40").Value2 = 2277.2222
$ws.Range("M40").Value2 = -1295.4117
$ws.Range("N40").Value2 = -2627.2222

$ws.Range("H51").Value2 = 3800
$ws.Range("I51").Value2 = 3800
$ws.Range("K51").Value2 = 3800
$ws.Range("M51").Value2 = -3316

$ws.Range("H62").Value2 = 2042.5
$ws.Range("I62").Value2 = 0
$ws.Range("K62").Value2 = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value2 = 2042.5
$ws.Range("I65").Value2 = 0
$ws.Range("K65").Value2 = 0
$ws.Range("M65").ClearContents()

$ws.Range("H70").Value2 = 1750
$ws.Range("I70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value2 = 1750
$ws.Range("I73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("M73").ClearContents()

$ws.Range("H104").Value2 = 218
$ws.Range("I104").Value2 = 227
$ws.Range("K104").Value2 = 681
$ws.Range("M104").Value2 = 1066

$ws.Range("H132").Value2 = 4333
$ws.Range("I132").Value2 = 4333
$ws.Range("K132").Value2 = 12999
$ws.Range("M132").Value2 = -10469

$ws.Range("H140").Value2 = 55000
$ws.Range("J140").Value2 = 55000
$ws.Range("L140").Value2 = 55000
$ws.Range("N140").Value2 = -65360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value2 = 2875
$ws.Range("I36").Value2 = 2875
$ws.Range("K36").Value2 = 2875
$ws.Range("M36").Value2 = -2529

$ws.Range("H45").Value2 = 1347.25
$ws.Range("I45").Value2 = 1347.25
$ws.Range("J45").Value2 = 0
$ws.Range("K45").Value2 = 1347.25
$ws.Range("L45").Value2 = 0
$ws.Range("M45").Value2 = -970.25
$ws.Range("N45").ClearContents()

$ws.Range("H97").Value2 = 611.8461
$ws.Range("J97").Value2 = 531.8333
$ws.Range("L97").Value2 = 531.8333
$ws.Range("N97").Value2 = -1523.8333

$ws.Range("H110").Value2 = 0
$ws.Range("I110").Value2 = 0
$ws.Range("K110").Value2 = 0
$ws.Range("M110").ClearContents()

$ws.Range("H139").Value2 = 69890
$ws.Range("J139").Value2 = 69890
$ws.Range("L139").Value2 = 69890
$ws.Range("N139").Value2 = -80170

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 205.71428
$ws.Range("I22").Value2 = 156.66667
$ws.Range("K22").Value2 = 156.66667
$ws.Range("M22").Value2 = 16.33332999999999

$ws.Range("H36").Value2 = 3050
$ws.Range("I36").Value2 = 1111
$ws.Range("J36").Value2 = 4989
$ws.Range("K36").Value2 = 1111
$ws.Range("L36").Value2 = 4989
$ws.Range("M36").Value2 = -577
$ws.Range("N36").Value2 = -6057

$ws.Range("H99").Value2 = 2374.875
$ws.Range("I99").Value2 = 2428.5715
$ws.Range("K99").Value2 = 2428.5715
$ws.Range("M99").Value2 = -930.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 398.5
$ws.Range("I22").Value2 = 398.5
$ws.Range("K22").Value2 = 398.5
$ws.Range("M22").Value2 = -48.5

$ws.Range("H41").Value2 = 999.3333
$ws.Range("I41").Value2 = 999.3333
$ws.Range("K41").Value2 = 999.3333
$ws.Range("M41").Value2 = -571.3333

$ws.Range("H42").Value2 = 8499.666999999999
$ws.Range("I42").Value2 = 4249.5
$ws.Range("K42").Value2 = 4249.5
$ws.Range("M42").Value2 = -3656.5

$ws.Range("H60").Value2 = 25000
$ws.Range("J60").Value2 = 25000
$ws.Range("L60").Value2 = 25000
$ws.Range("N60").Value2 = -26022

$ws.Range("H62").Value2 = 4033
$ws.Range("I62").Value2 = 3799
$ws.Range("K62").Value2 = 3799
$ws.Range("M62").Value2 = -3175

$ws.Range("H65").Value2 = 4033
$ws.Range("I65").Value2 = 3799
$ws.Range("K65").Value2 = 18995
$ws.Range("M65").Value2 = -15875

$ws.Range("H105").Value2 = 2409.35
$ws.Range("I105").Value2 = 1511.3
$ws.Range("K105").Value2 = 1511.3
$ws.Range("M105").Value2 = 235.7

$ws.Range("H118").Value2 = 57000
$ws.Range("J118").Value2 = 57000
$ws.Range("L118").Value2 = 57000
$ws.Range("N118").Value2 = -60314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value2 = 2000
$ws.Range("J64").Value2 = 2000
$ws.Range("L64").Value2 = 6000
$ws.Range("N64").Value2 = -6540

$ws.Range("H67").Value2 = 2000
$ws.Range("J67").Value2 = 2000
$ws.Range("L67").Value2 = 6000
$ws.Range("N67").Value2 = -7872

$ws.Range("H111").Value2 = 2555
$ws.Range("I111").Value2 = 2555
$ws.Range("K111").Value2 = 7665
$ws.Range("M111").Value2 = -4598

$ws.Range("H113").Value2 = 1237.2858
$ws.Range("J113").Value2 = 1276
$ws.Range("L113").Value2 = 3828
$ws.Range("N113").Value2 = -8168

$ws.Range("H138").Value2 = 1597.8
$ws.Range("I138").Value2 = 1597.8
$ws.Range("K138").Value2 = 4793.4
$ws.Range("M138").Value2 = 346.6000000000004

$ws.Range("H139").Value2 = 2914.375
$ws.Range("I139").Value2 = 1132.3334
$ws.Range("J139").Value2 = 3983.6
$ws.Range("K139").Value2 = 3397.0002
$ws.Range("L139").Value2 = 11950.8
$ws.Range("M139").Value2 = 1742.9998
$ws.Range("N139").Value2 = -22230.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 124.36842
$ws.Range("J2").Value2 = 121.333336
$ws.Range("L2").Value2 = 121.333336
$ws.Range("N2").Value2 = -347.333336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value2 = 898.3333
$ws.Range("I68").Value2 = 1246
$ws.Range("J68").Value2 = 724.5
$ws.Range("K68").Value2 = 1246
$ws.Range("L68").Value2 = 724.5
$ws.Range("M68").Value2 = -497
$ws.Range("N68").Value2 = -2222.5

$ws.Range("H71").Value2 = 898.3333
$ws.Range("I71").Value2 = 1246
$ws.Range("J71").Value2 = 724.5
$ws.Range("K71").Value2 = 6230
$ws.Range("L71").Value2 = 3622.5
$ws.Range("M71").Value2 = -2486
$ws.Range("N71").Value2 = -11110.5

$ws.Range("H82").Value2 = 1624.5
$ws.Range("I82").Value2 = 1724.5
$ws.Range("J82").Value2 = 1524.5
$ws.Range("K82").Value2 = 1724.5
$ws.Range("L82").Value2 = 1524.5
$ws.Range("M82").Value2 = -1363.5
$ws.Range("N82").Value2 = -2246.5

$ws.Range("H85").Value2 = 1624.5
$ws.Range("I85").Value2 = 1724.5
$ws.Range("J85").Value2 = 1524.5
$ws.Range("K85").Value2 = 1724.5
$ws.Range("L85").Value2 = 1524.5
$ws.Range("M85").Value2 = -476.5
$ws.Range("N85").Value2 = -4020.5

$ws.Range("H100").Value2 = 3228.3572
$ws.Range("I100").Value2 = 3154.4546
$ws.Range("K100").Value2 = 3154.4546
$ws.Range("M100").Value2 = -2613.4546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value2 = 4494.75
$ws.Range("J30").Value2 = 4990
$ws.Range("L30").Value2 = 4990
$ws.Range("N30").Value2 = -5204

$ws.Range("H41").Value2 = 19700.334
$ws.Range("I41").Value2 = 19676
$ws.Range("K41").Value2 = 19676
$ws.Range("M41").Value2 = -19286

$ws.Range("H136").Value2 = 2241.5
$ws.Range("I136").Value2 = 2638.8
$ws.Range("J136").Value2 = 1248.25
$ws.Range("K136").Value2 = 7916.400000000001
$ws.Range("L136").Value2 = 3744.75
$ws.Range("M136").Value2 = -5366.400000000001
$ws.Range("N136").Value2 = -8844.75
